# "Changed day 3 topic order"
#
# Slide 1 (title slide, Rectangle 2): "Tag 3: GitOps,<br/>Docker, Deployment-Strategien"
#   -> "Tag 3: Docker, GitOps, Deployment-Strategien"  (drop the manual line break,
#      move "Docker, " in front of "GitOps")
#
# Slide 2 & 3 (Agenda, "Inhaltsplatzhalter 18"): the "Tag 3 - ..." bullet and its
#   four sub-bullets get the same reshuffle:
#     "Tag 3 - GitOps, Docker, Deployment-Strategien"  -> "Tag 3 - Docker, GitOps, Deployment-Strategien"
#     "GitOps Grundlagen"                              -> "Entwicklung mit Docker"
#     "Entwicklung mit Docker"                         -> "Container/Docker-Registry"
#     "Container/Docker-Registry"                      -> "Erstellen von Release- und Tagged-Images"
#     "Erstellen von Release- und Tagged-Images"       -> "GitOps Grundlagen"

$dash = [char]0x2013

function Update-Tag3Slide($slideIndex) {
    $s  = $ppt.ActivePresentation.Slides.Item($slideIndex)
    $sh = $s.Shapes.Item(1)
    $tr = $sh.TextFrame.TextRange

    # --- Paragraph 13: "Tag 3 - GitOps, Docker, Deployment-Strategien" ---------
    # Runs: "Tag 3 - " | "GitOps"(err) | ", Docker, " | "Deployment"(err) | "-Strategien"
    # Move "Docker, " from the 3rd run to the end of the 1st run, leaving ", " behind.
    # (processed right-to-left so earlier offsets stay valid)
    $tr.Characters(333, 10).Text = ", "
    $tr.Characters(319, 8).Text  = "Tag 3 " + $dash + " Docker, "

    # --- Paragraphs 14-17: rotate the four sub-bullets -------------------------
    # Processed from the last paragraph back to the first so offsets of the
    # paragraphs not yet touched stay valid.

    # Paragraph 17: "Erstellen von Release- und " + "Tagged"(err) + "-Images"
    #            -> "GitOps"(err) + " Grundlagen"
    # Drop the leading "Erstellen von Release- und " run entirely, then reuse the
    # remaining two runs (which already carry the right err/no-err split) in place.
    $tr.Characters(432, 27).Text = ""
    $tr.Characters(438, 7).Text  = " Grundlagen"
    $tr.Characters(432, 6).Text  = "GitOps"

    # Paragraph 16: "Container/Docker-Registry" -> "Erstellen von Release- und Tagged-Images"
    $tr.Characters(406, 25).Text = "Erstellen von Release- und Tagged-Images"

    # Paragraph 15: "Entwicklung mit Docker" -> "Container/Docker-Registry"
    $tr.Characters(383, 22).Text = "Container/Docker-Registry"

    # Paragraph 14: "GitOps"(err) + " Grundlagen" -> "Entwicklung mit Docker"
    # Drop the leading "GitOps" run (it carries the err flag we don't want here),
    # then reuse the remaining " Grundlagen" run (no err) for the new text.
    $tr.Characters(365, 6).Text  = ""
    $tr.Characters(365, 11).Text = "Entwicklung mit Docker"
}

# --- Slide 1: title slide --------------------------------------------------
$s1  = $ppt.ActivePresentation.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange
$tr1.Characters(1, 44).Text = "Tag 3: Docker, GitOps, Deployment-Strategien"

# --- Slides 2 & 3: Agenda overview -----------------------------------------
Update-Tag3Slide 2
Update-Tag3Slide 3
